# Final update to lab 4 with related source code files
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "Students"

# Update individual score cells
$ws.Range("D2").Value = 64

$ws.Range("B3").Value = 58
$ws.Range("C3").Value = 65
$ws.Range("D3").Value = 85
$ws.Range("E3").Value = 38

$ws.Range("C4").Value = 62
$ws.Range("D4").Value = 58

$ws.Range("D5").Value = 55
$ws.Range("E5").Value = 61

$ws.Range("B6").Value = 65
$ws.Range("C6").Value = 54

$ws.Range("C10").Value = 69
$ws.Range("D10").Value = 57

$ws.Range("B11").Value = 35
$ws.Range("C11").Value = 80
$ws.Range("D11").Value = 65
$ws.Range("E11").Value = 58

# Remove the last 5 student rows (Student11 - Student15)
$ws.Range("A12:E16").EntireRow.Delete()

# Update the view: zoom + active selection
$excel.ActiveWindow.Zoom = 205
$null = $ws.Range("A13").Select()
